{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Applies the edit described by the diff:\n//  1. \"Titre : ... Brumadinho ...\" -> wrap \"Brumadinho\" with proofErr spellStart/spellEnd\n//  2. \"Emilien Valin; Macosso Michael; ...\" -> wrap \"Macosso\" with proofErr spellStart/spellEnd\n//  3. \"1)\" paragraph -> becomes \"2. Introduction\" + a large new block of inserted\n//     paragraphs (Introduction + Fiche d'identit\u00e9 tableau synth\u00e9tique)\n//  4. Move <w:lastRenderedPageBreak/> from the \"r\u00e9diger des recommandations...\"\n//     paragraph to the \"telles catastrophes...\" paragraph.\n//  5. \"Je reviens du monde d'avant france inter - Serie de reportages de Giv Anquetil.\"\n//     -> wrap \"france\", \"Serie\", \"Giv\" with proofErr spellStart/spellEnd\n//  6. \"Site de FGH Sciences Humaines ... Brumadinho.\" -> wrap \"Brumadinho\" with proofErr\n//  7. \"- Times ou calibri 12\" -> wrap \"calibri\" with proofErr spellStart/spellEnd\n//\n// Because several of these edits change the total paragraph count (the big \"1)\"\n// block especially), paragraphs are re-located by their distinctive text right\n// before each mutation instead of relying on indices captured up front.\n\n// Helper: wrap one or more raw <w:p>...</w:p> fragments (WordprocessingML, using the\n// \"w\" namespace prefix) into the Flat OPC package format that insertOoxml expects.\nfunction flatOpc(bodyFragment) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyFragment + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst body = context.document.body;\n\n// Find (and return) the paragraph whose text satisfies `predicate`. Re-queries the\n// live paragraph collection every time it is called so indices shifted by earlier\n// edits never go stale.\nasync function findParagraph(predicate) {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (predicate(paragraphs.items[i].text)) return paragraphs.items[i];\n  }\n  throw new Error(\"paragraph not found\");\n}\n\n// 1) Titre paragraph: split \"Brumadinho\" out with proofErr markers.\nlet p = await findParagraph((t) => t.indexOf(\"Titre : Rupture du barrage de Brumadinho\") !== -1);\np.insertOoxml(\n  flatOpc(\n    '<w:p><w:r><w:t xml:space=\"preserve\">Titre : Rupture du barrage de </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> (2019) : un al\u00e9a technologique majeur</w:t></w:r></w:p>'\n  ),\n  \"Replace\"\n);\nawait context.sync();\n\n// 2) Emilien Valin paragraph: split \"Macosso\" out with proofErr markers (rest of the\n//    paragraph, i.e. the \"suppl\u00e9mentaires\" runs, is untouched).\np = await findParagraph((t) => t.indexOf(\"Emilien Valin; Macosso Michael;\") !== -1);\np.insertOoxml(\n  flatOpc(\n    '<w:p><w:r><w:t xml:space=\"preserve\">Emilien Valin; </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Macosso</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> Michael; (deux types </w:t></w:r>' +\n      '<w:r><w:t>suppl\u00e9mentaires</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> si ils foutent le nez ici et mettent leurs noms)</w:t></w:r></w:p>'\n  ),\n  \"Replace\"\n);\nawait context.sync();\n\n// 3) \"1)\" paragraph becomes \"2. Introduction\" immediately followed by the new\n//    Introduction + Fiche d'identit\u00e9 (tableau synth\u00e9tique) content block.\nconst introBlock =\n  '<w:p><w:r><w:t>2. Introduction</w:t></w:r></w:p>' +\n  \"<w:p/>\" +\n  '<w:p><w:r><w:t>D\u00e9finition al\u00e9a technologique : danger li\u00e9 \u00e0 une activit\u00e9 humaine pouvant causer des dommages majeurs.</w:t></w:r></w:p>' +\n  \"<w:p/>\" +\n  '<w:p><w:r><w:t xml:space=\"preserve\">Exemple : rupture d\u2019un barrage de r\u00e9sidus miniers \u00e0 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>, Br\u00e9sil, 25 janvier 2019.</w:t></w:r></w:p>' +\n  \"<w:p/><w:p/>\" +\n  '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>' +\n  \"<w:p/><w:p/><w:p/>\" +\n  '<w:p><w:r><w:t>3. Fiche d\u2019identit\u00e9 (tableau synth\u00e9tique)</w:t></w:r></w:p>' +\n  \"<w:p/>\" +\n  '<w:p><w:r><w:t>\u00c9v\u00e9nement</w:t></w:r>' +\n  '<w:r><w:tab/><w:t xml:space=\"preserve\">Rupture du barrage de r\u00e9sidus (Dam B1) \u2013 </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>' +\n  \"<w:p/>\" +\n  '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Date</w:t></w:r><w:r><w:tab/><w:t>25 janvier 2019</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Lieu</w:t></w:r><w:r><w:tab/><w:t xml:space=\"preserve\">Mine </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>C\u00f3rrego</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> do </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Feij\u00e3o</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>, Minas Gerais, Br\u00e9sil</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Op\u00e9rateur</w:t></w:r><w:r><w:tab/><w:t>Vale S.A.</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Volume de boues d\u00e9vers\u00e9es</w:t></w:r><w:r><w:tab/><w:t>~11\u201312 millions m\u00b3</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Victimes</w:t></w:r><w:r><w:tab/><w:t>Environ 270 morts</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Cause technique</w:t></w:r><w:r><w:tab/><w:t>Liqu\u00e9faction statique des r\u00e9sidus (</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>static</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>liquefaction</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>)</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Impacts principaux</w:t></w:r><w:r><w:tab/>' +\n  '<w:t>Pertes humaines, destruction d\u2019infrastructures et habitations, contamination des cours d\u2019eau, impacts socio-\u00e9conomiques</w:t></w:r></w:p>' +\n  \"<w:p/><w:p/>\" +\n  '<w:p><w:r><w:t xml:space=\"preserve\">Sources : </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Wikipedia</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u00ab </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> dam </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>disaster</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> \u00bb, Global Tailings Portal, rapports techniques post-accident.</w:t></w:r></w:p>' +\n  \"<w:p/>\" +\n  '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>';\n\np = await findParagraph((t) => t.trim() === \"1)\");\np.insertOoxml(flatOpc(introBlock), \"Replace\");\nawait context.sync();\n\n// 4) Move the page-break marker: remove it from \"r\u00e9diger...\" and add it to the\n//    start of \"telles catastrophes...\".\np = await findParagraph((t) => t.indexOf(\"r\u00e9diger des recommandations aux autorit\u00e9s\") !== -1);\np.insertOoxml(\n  flatOpc('<w:p><w:r><w:t>r\u00e9diger des recommandations aux autorit\u00e9s du pays concern\u00e9 afin d\u2019\u00e9viter de</w:t></w:r></w:p>'),\n  \"Replace\"\n);\nawait context.sync();\n\np = await findParagraph((t) => t.indexOf(\"telles catastrophes \u00e0 l\u2019avenir\") !== -1);\np.insertOoxml(\n  flatOpc('<w:p><w:r><w:lastRenderedPageBreak/><w:t>telles catastrophes \u00e0 l\u2019avenir. Minimum 15 lignes.</w:t></w:r></w:p>'),\n  \"Replace\"\n);\nawait context.sync();\n\n// 5) \"Je reviens du monde d'avant france inter - Serie de reportages de Giv Anquetil.\"\np = await findParagraph((t) => t.indexOf(\"Je reviens du monde d'avant\") !== -1);\np.insertOoxml(\n  flatOpc(\n    '<w:p><w:r><w:t xml:space=\"preserve\">Je reviens du monde d\\'avant </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>france</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> inter</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> - </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Serie</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> de reportages de </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Giv</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> Anquetil</w:t></w:r>' +\n      '<w:r><w:t>.</w:t></w:r></w:p>'\n  ),\n  \"Replace\"\n);\nawait context.sync();\n\n// 6) \"Site de FGH Sciences Humaines ... li\u00e9 aux mines de Brumadinho.\"\np = await findParagraph((t) => t.indexOf(\"Site de FGH Sciences Humaines\") !== -1);\np.insertOoxml(\n  flatOpc(\n    '<w:p><w:r><w:t xml:space=\"preserve\">Site de FGH Sciences Humaines disposant de carte et d\u2019un dossier th\u00e9orique sur la mise en \u00e9vidence et la vulgarisation de m\u00e9thode de travail cartographique sur le sujet des risque et al\u00e9a technologique li\u00e9 aux mines de </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t>.</w:t></w:r></w:p>'\n  ),\n  \"Replace\"\n);\nawait context.sync();\n\n// 7) \"- Times ou calibri 12\"\np = await findParagraph((t) => t.trim() === \"- Times ou calibri 12\");\np.insertOoxml(\n  flatOpc(\n    '<w:p><w:r><w:t xml:space=\"preserve\">- Times ou </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>calibri</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> 12</w:t></w:r></w:p>'\n  ),\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument / $d is the live document.\n#\n# Applies the edit described by the diff:\n#  1. \"Titre : ... Brumadinho ...\" -> wrap \"Brumadinho\" with proofErr spellStart/spellEnd\n#  2. \"Emilien Valin; Macosso Michael; ...\" -> wrap \"Macosso\" with proofErr spellStart/spellEnd\n#  3. \"1)\" paragraph -> becomes \"2. Introduction\" + a large new block of inserted\n#     paragraphs (Introduction + Fiche d'identit\u00e9 tableau synth\u00e9tique)\n#  4. Move <w:lastRenderedPageBreak/> from the \"r\u00e9diger des recommandations...\"\n#     paragraph to the \"telles catastrophes...\" paragraph.\n#  5. \"Je reviens du monde d'avant france inter - Serie de reportages de Giv Anquetil.\"\n#     -> wrap \"france\", \"Serie\", \"Giv\" with proofErr spellStart/spellEnd\n#  6. \"Site de FGH Sciences Humaines ... Brumadinho.\" -> wrap \"Brumadinho\" with proofErr\n#  7. \"- Times ou calibri 12\" -> wrap \"calibri\" with proofErr spellStart/spellEnd\n#\n# Paragraphs are re-located by their distinctive text right before each mutation\n# (instead of a fixed index) since edits such as the big \"1)\" block change the\n# total paragraph count and would otherwise shift later indices.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParaByText($doc, $substr) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text -like \"*$substr*\") {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Wrap-FlatOpc($bodyFragment) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# 1) Titre paragraph: split \"Brumadinho\" out with proofErr markers.\n$p = Find-ParaByText $d \"Titre : Rupture du barrage de Brumadinho\"\n$frag = '<w:p><w:r><w:t xml:space=\"preserve\">Titre : Rupture du barrage de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> (2019) : un al\u00e9a technologique majeur</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n# 2) Emilien Valin paragraph: split \"Macosso\" out with proofErr markers (rest of the\n#    paragraph, i.e. the \"suppl\u00e9mentaires\" runs, is untouched).\n$p = Find-ParaByText $d \"Emilien Valin; Macosso Michael;\"\n$frag = '<w:p><w:r><w:t xml:space=\"preserve\">Emilien Valin; </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Macosso</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> Michael; (deux types </w:t></w:r><w:r><w:t>suppl\u00e9mentaires</w:t></w:r><w:r><w:t xml:space=\"preserve\"> si ils foutent le nez ici et mettent leurs noms)</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n# 3) \"1)\" paragraph becomes \"2. Introduction\" immediately followed by the new\n#    Introduction + Fiche d'identit\u00e9 (tableau synth\u00e9tique) content block.\n$introBlock = '<w:p><w:r><w:t>2. Introduction</w:t></w:r></w:p>'\n$introBlock += '<w:p/>'\n$introBlock += '<w:p><w:r><w:t>D\u00e9finition al\u00e9a technologique : danger li\u00e9 \u00e0 une activit\u00e9 humaine pouvant causer des dommages majeurs.</w:t></w:r></w:p>'\n$introBlock += '<w:p/>'\n$introBlock += '<w:p><w:r><w:t xml:space=\"preserve\">Exemple : rupture d\u2019un barrage de r\u00e9sidus miniers \u00e0 </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>, Br\u00e9sil, 25 janvier 2019.</w:t></w:r></w:p>'\n$introBlock += '<w:p/><w:p/>'\n$introBlock += '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>'\n$introBlock += '<w:p/><w:p/><w:p/>'\n$introBlock += '<w:p><w:r><w:t>3. Fiche d\u2019identit\u00e9 (tableau synth\u00e9tique)</w:t></w:r></w:p>'\n$introBlock += '<w:p/>'\n$introBlock += '<w:p><w:r><w:t>\u00c9v\u00e9nement</w:t></w:r><w:r><w:tab/><w:t xml:space=\"preserve\">Rupture du barrage de r\u00e9sidus (Dam B1) \u2013 </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n$introBlock += '<w:p/>'\n$introBlock += '<w:p><w:r><w:lastRenderedPageBreak/><w:t>Date</w:t></w:r><w:r><w:tab/><w:t>25 janvier 2019</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Lieu</w:t></w:r><w:r><w:tab/><w:t xml:space=\"preserve\">Mine </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>C\u00f3rrego</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> do </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Feij\u00e3o</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>, Minas Gerais, Br\u00e9sil</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Op\u00e9rateur</w:t></w:r><w:r><w:tab/><w:t>Vale S.A.</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Volume de boues d\u00e9vers\u00e9es</w:t></w:r><w:r><w:tab/><w:t>~11\u201312 millions m\u00b3</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Victimes</w:t></w:r><w:r><w:tab/><w:t>Environ 270 morts</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Cause technique</w:t></w:r><w:r><w:tab/><w:t>Liqu\u00e9faction statique des r\u00e9sidus (</w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>static</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>liquefaction</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>)</w:t></w:r></w:p>'\n$introBlock += '<w:p><w:r><w:t>Impacts principaux</w:t></w:r><w:r><w:tab/><w:t>Pertes humaines, destruction d\u2019infrastructures et habitations, contamination des cours d\u2019eau, impacts socio-\u00e9conomiques</w:t></w:r></w:p>'\n$introBlock += '<w:p/><w:p/>'\n$introBlock += '<w:p><w:r><w:t xml:space=\"preserve\">Sources : </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Wikipedia</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> \u00ab </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> dam </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>disaster</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> \u00bb, Global Tailings Portal, rapports techniques post-accident.</w:t></w:r></w:p>'\n$introBlock += '<w:p/>'\n$introBlock += '<w:p><w:pPr><w:pBdr><w:bottom w:val=\"single\" w:sz=\"6\" w:space=\"1\" w:color=\"auto\"/></w:pBdr></w:pPr></w:p>'\n\n$p = Find-ParaByText $d \"1)\"\n$p.Range.InsertXML((Wrap-FlatOpc $introBlock))\n\n# 4) Move the page-break marker: remove it from \"r\u00e9diger...\" and add it to the\n#    start of \"telles catastrophes...\".\n$p = Find-ParaByText $d \"r\u00e9diger des recommandations aux autorit\u00e9s\"\n$frag = '<w:p><w:r><w:t>r\u00e9diger des recommandations aux autorit\u00e9s du pays concern\u00e9 afin d\u2019\u00e9viter de</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n$p = Find-ParaByText $d \"telles catastrophes \u00e0 l\u2019avenir\"\n$frag = '<w:p><w:r><w:lastRenderedPageBreak/><w:t>telles catastrophes \u00e0 l\u2019avenir. Minimum 15 lignes.</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n# 5) \"Je reviens du monde d'avant france inter - Serie de reportages de Giv Anquetil.\"\n$p = Find-ParaByText $d \"Je reviens du monde d'avant\"\n$frag = '<w:p><w:r><w:t xml:space=\"preserve\">Je reviens du monde d''avant </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>france</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> inter</w:t></w:r><w:r><w:t xml:space=\"preserve\"> - </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Serie</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> de reportages de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Giv</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> Anquetil</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n# 6) \"Site de FGH Sciences Humaines ... li\u00e9 aux mines de Brumadinho.\"\n$p = Find-ParaByText $d \"Site de FGH Sciences Humaines\"\n$frag = '<w:p><w:r><w:t xml:space=\"preserve\">Site de FGH Sciences Humaines disposant de carte et d\u2019un dossier th\u00e9orique sur la mise en \u00e9vidence et la vulgarisation de m\u00e9thode de travail cartographique sur le sujet des risque et al\u00e9a technologique li\u00e9 aux mines de </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Brumadinho</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>.</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n\n# 7) \"- Times ou calibri 12\"\n$p = Find-ParaByText $d \"- Times ou calibri 12\"\n$frag = '<w:p><w:r><w:t xml:space=\"preserve\">- Times ou </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>calibri</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> 12</w:t></w:r></w:p>'\n$p.Range.InsertXML((Wrap-FlatOpc $frag))\n"}
